# Update "想去人数" (people interested) counters after a fresh scrape,
# as reflected in the gh-pages output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5935
$wsExhibit.Range("F4").Value = 1143
$wsExhibit.Range("F11").Value = 37
$wsExhibit.Range("F13").Value = 2117
$wsExhibit.Range("F14").Value = 1534
$wsExhibit.Range("F15").Value = 1184
$wsExhibit.Range("F18").Value = 460
$wsExhibit.Range("F21").Value = 1077
$wsExhibit.Range("F24").Value = 3895
$wsExhibit.Range("F29").Value = 59
$wsExhibit.Range("F30").Value = 565
$wsExhibit.Range("F36").Value = 876
$wsExhibit.Range("F39").Value = 95

# Sheet "演出" (rId2 / sheet2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 763

# Sheet "全部类型" (rId4 / sheet4) - combined view of all events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5935
$wsAll.Range("F4").Value = 1143
$wsAll.Range("F6").Value = 763
$wsAll.Range("F15").Value = 37
$wsAll.Range("F18").Value = 2117
$wsAll.Range("F19").Value = 1534
$wsAll.Range("F20").Value = 1184
$wsAll.Range("F23").Value = 460
$wsAll.Range("F27").Value = 1077
$wsAll.Range("F30").Value = 3895
$wsAll.Range("F35").Value = 59
$wsAll.Range("F36").Value = 565
$wsAll.Range("F42").Value = 876
$wsAll.Range("F45").Value = 95
